# feat: generate Excel report for reorder strategy recommendations
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): bold, centered/top-aligned, thin box border ---
$ws.Range("A1").Value = "category"
$ws.Range("B1").Value = "recommendation"
$ws.Range("C1").Value = "new_safety_stock"
$ws.Range("D1").Value = "new_reorder_point"
$ws.Range("E1").Value = "new_optimal_inventory"
$ws.Range("F1").Value = "new_holding_cost"
$ws.Range("G1").Value = "potential_saving"

$headerRange = $ws.Range("A1:G1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# --- Data rows ---
$ws.Range("A2").Value = "toys"
$ws.Range("B2").Value = "Giảm Safety Stock từ 93263 → 74610 và Reorder Point từ 507836 → 457052 để tiết kiệm chi phí."
$ws.Range("C2").Value = 74610
$ws.Range("D2").Value = 457052
$ws.Range("E2").Value = 531662
$ws.Range("F2").Value = 1169656
$ws.Range("G2").Value = -640689

$ws.Range("A3").Value = "garden_tools"
$ws.Range("B3").Value = "Giảm Safety Stock từ 1552 → 1241 và Reorder Point từ 189523 → 170570 để tiết kiệm chi phí."
$ws.Range("C3").Value = 1241
$ws.Range("D3").Value = 170570
$ws.Range("E3").Value = 171811
$ws.Range("F3").Value = 412346
$ws.Range("G3").Value = -228914
